# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the freshly generated data snapshot (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 13664
$wsExhibit.Range("F9").Value = 13873
$wsExhibit.Range("F10").Value = 14688

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 13664
$wsAll.Range("F10").Value = 13873
$wsAll.Range("F11").Value = 14688
